# Update the Burundi MSME summary numbers in row 11 ("Enterprises density
# (per 1000 people)") and row 12 ("Enterprises (% of total)").
#
# The source values are stored as literal text (shared strings) rather than
# numbers, so a plain `.Value = "0.37"` assignment would let Excel's normal
# type-inference turn them into numeric cells. To keep them as text without
# leaving a stray NumberFormat/quote-prefix style on the cell, each cell is
# briefly switched to the Text format, given its new literal value, and then
# restored to the workbook's "Normal" cell style (General format / default
# styling) that it started with.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 11: Enterprises density (per 1000 people)
Set-TextValue "B11" "0.37"
Set-TextValue "C11" "0.04"
Set-TextValue "D11" "0.41"

# Row 12: Enterprises (% of total)
Set-TextValue "B12" "89.01"
Set-TextValue "C12" "10.13"
Set-TextValue "D12" "99.14"
